# Natmi following Dr Hou advice
# Rebuilds the Efna2->Epha7 sending/target cluster results table (rows 2-10,
# columns A:T) with the updated values, including the new sCs-sourced rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("ECs","Efna2","Epha7","ECs",3,1,1.408030333333333,4.224091,0.3454737251382253,0.3454737251382253,2,0.6666666666666666,0.1532026666666667,0.459608,0.03935316038222934,0.03935316038222934,0.2157140018142222,1.941426016328,0.01359548291321079,0.0135954829132108)
  ,@("ECs","Efna2","Epha7","FAPs",3,1,1.408030333333333,4.224091,0.3454737251382253,0.3454737251382253,3,1,0.1693566666666667,0.50807,0.04350263745496,0.04350263745496,0.2384593238188889,2.14613391437,0.01502901821490272,0.01502901821490272)
  ,@("ECs","Efna2","Epha7","sCs",3,1,1.408030333333333,4.224091,0.3454737251382253,0.3454737251382253,3,1,3.570461333333333,10.711384,0.9171442021628107,0.9171442021628107,5.02731786132711,45.24586075194399,0.3168492240101118,0.3168492240101118)
  ,@("FAPs","Efna2","Epha7","ECs",3,1,2.015377,6.046131,0.494492045565236,0.4944920455652361,2,0.6666666666666666,0.1532026666666667,0.459608,0.03935316038222934,0.03935316038222934,0.3087611307386667,2.778850176648,0.01945982477686539,0.01945982477686539)
  ,@("FAPs","Efna2","Epha7","FAPs",3,1,2.015377,6.046131,0.494492045565236,0.4944920455652361,3,1,0.1693566666666667,0.50807,0.04350263745496,0.04350263745496,0.3413175307966667,3.07185777717,0.02151170818258602,0.02151170818258603)
  ,@("FAPs","Efna2","Epha7","sCs",3,1,2.015377,6.046131,0.494492045565236,0.4944920455652361,3,1,3.570461333333333,10.711384,0.9171442021628107,0.9171442021628107,7.195825650589332,64.76243085530399,0.4535205126057846,0.4535205126057846)
  ,@("sCs","Efna2","Epha7","ECs",3,1,0.6522436666666667,1.956731,0.1600342292965385,0.1600342292965385,2,0.6666666666666666,0.1532026666666667,0.459608,0.03935316038222934,0.03935316038222934,0.09992546904977778,0.899329221448,0.006297852692153146,0.006297852692153146)
  ,@("sCs","Efna2","Epha7","FAPs",3,1,0.6522436666666667,1.956731,0.1600342292965385,0.1600342292965385,3,1,0.1693566666666667,0.50807,0.04350263745496,0.04350263745496,0.1104618132411111,0.99415631917,0.006961911057471255,0.006961911057471255)
  ,@("sCs","Efna2","Epha7","sCs",3,1,0.6522436666666667,1.956731,0.1600342292965385,0.1600342292965385,3,1,3.570461333333333,10.711384,0.9171442021628107,0.9171442021628107,2.328810791744889,20.959297125704,0.1467744655469141,0.1467744655469141)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}
